# Applies the "Update post last call, with all examples and images" edit:
#  - Resize/reposition 3 connector arrows that point into the
#    SpecimenDefinition/PlanDefinition diagram boxes.
#  - Rename a few fields in the PlanDefinition pseudo-code text block:
#      specimenRequirements -> specimenRequested
#      sampleRequirement    -> material  (and tighten one brace's spacing)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    foreach ($i in 1..$slide.Shapes.Count) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) { return $shp }
    }
    return $null
}

# Shape.Left/Top/Width/Height are marshalled through a 32-bit float, so a
# naive "emu / 12700.0" assignment can land one EMU below the intended
# integer value once PowerPoint converts back to EMU on save. Nudge the
# point value up in very small increments until the value we read back
# truly round-trips to the exact target EMU before moving on.
function Set-ShapeEmuProperty($shape, $propName, $targetEmu) {
    $base = $targetEmu / 12700.0
    for ($k = 0; $k -le 50; $k++) {
        $cand = $base + ($k * 0.000001)
        switch ($propName) {
            "Left"   { $shape.Left = $cand }
            "Top"    { $shape.Top = $cand }
            "Width"  { $shape.Width = $cand }
            "Height" { $shape.Height = $cand }
        }
        switch ($propName) {
            "Left"   { $rb = $shape.Left }
            "Top"    { $rb = $shape.Top }
            "Width"  { $rb = $shape.Width }
            "Height" { $rb = $shape.Height }
        }
        $emuBack = [Math]::Floor(($rb * 12700.0) + 0.5)
        if ([Math]::Abs($emuBack - $targetEmu) -lt 1) {
            return
        }
    }
}

function Set-ConnectorGeometry($slide, $shapeId, $left, $top, $width, $height) {
    $shp = Get-ShapeById $slide $shapeId
    Set-ShapeEmuProperty $shp "Left" $left
    Set-ShapeEmuProperty $shp "Top" $top
    Set-ShapeEmuProperty $shp "Width" $width
    Set-ShapeEmuProperty $shp "Height" $height
}

# Connecteur droit avec flèche 33 (id 34) -> "24h urine" box
Set-ConnectorGeometry $s 34 5045901 2155977 1502371 1774377
# Connecteur droit avec flèche 44 (id 45) -> "creatinine clearance" box
Set-ConnectorGeometry $s 45 5121965 718416 1688694 2004906
# Connecteur droit avec flèche 77 (id 78) -> "title/custodian" box
Set-ConnectorGeometry $s 78 5121965 1433332 1688694 1524638

# --- Text updates inside the PlanDefinition pseudo-code box (shape id 36) ---
function Replace-FirstOccurrence($textRange, $oldStr, $newStr) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldStr)
    if ($idx -lt 0) {
        return
    }
    $sub = $textRange.Characters($idx + 1, $oldStr.Length)
    $sub.Text = $newStr
}

$shp36 = Get-ShapeById $s 36
$tr = $shp36.TextFrame.TextRange

Replace-FirstOccurrence $tr "specimenRequirements [ " "specimenRequested [ "
Replace-FirstOccurrence $tr "             { sampleRequirement  [] }, " "             { material  [] }, "
Replace-FirstOccurrence $tr "             { sampleRequirement  [] }" "             {material  [] }"
Replace-FirstOccurrence $tr "             { sampleRequirement  [] }" "             {material  [] }"
